# Update "想去人数" (F column) figures and one Cover image URL (I column)
# on the "展览" (sheet1) and "全部类型" (sheet4) worksheets, reflecting a
# newer data pull (see commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$newCover = "//i2.hdslb.com/bfs/openplatform/202409/xp4jNVRG1727165677359.jpeg"

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value = 60
$ws1.Range("F5").Value = 693
$ws1.Range("F6").Value = 70
$ws1.Range("F7").Value = 2113
$ws1.Range("F10").Value = 4669
$ws1.Range("F16").Value = 153
$ws1.Range("F20").Value = 3564
$ws1.Range("F21").Value = 137
$ws1.Range("F25").Value = 92
$ws1.Range("F29").Value = 76
$ws1.Range("F30").Value = 214
$ws1.Range("F31").Value = 17
$ws1.Range("I31").Value = $newCover
$ws1.Range("F32").Value = 769
$ws1.Range("F33").Value = 2212
$ws1.Range("F34").Value = 410

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F4").Value = 60
$ws4.Range("F5").Value = 693
$ws4.Range("F6").Value = 70
$ws4.Range("F7").Value = 2113
$ws4.Range("F10").Value = 4669
$ws4.Range("F16").Value = 153
$ws4.Range("F20").Value = 3564
$ws4.Range("F21").Value = 137
$ws4.Range("F25").Value = 92
$ws4.Range("F29").Value = 76
$ws4.Range("F30").Value = 214
$ws4.Range("F32").Value = 17
$ws4.Range("I32").Value = $newCover
$ws4.Range("F33").Value = 769
$ws4.Range("F34").Value = 2212
$ws4.Range("F35").Value = 410
